$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 3 contact info value (was 8394870189, should now be 7078140054)
$ws.Range("B3").Value = 7078140054

# Delete rows 4 and 5 (duplicate "Devesh Rawat" entries)
$ws.Range("A4:B5").EntireRow.Delete()

# Update the active selection to F10 as shown in the diff
$ws.Range("F10").Select()
